# Auto generate OTPUB and error handling.
#
# - Strip the "(R)" / registered-trademark glyph from every FedEx sheet
#   name (Excel/Windows choked on the special character during the
#   automated OTPUB generation, so the names are being normalized to
#   plain ASCII).
# - Re-point the active sheet/tab at "Intra-Hawaii Standard List Rate"
#   (was previously left on "FHD Ground").

$wb = $excel.ActiveWorkbook

$renames = @{
    "FedEx First Overnight®"          = "FedEx First Overnight"
    "FedEx Priority Overnight®"       = "FedEx Priority Overnight"
    "FedEx Standard Overnight®"       = "FedEx Standard Overnight"
    "FedEx 2Day® A.M."                = "FedEx 2Day A.M."
    "FedEx 2Day®"                     = "FedEx 2Day"
    "FedEx Express Saver®"            = "FedEx Express Saver"
    "MW FedEx First Overnight®"       = "MW FedEx First Overnight"
    "MW FedEx Priority Overnight®"    = "MW FedEx Priority Overnight"
    "MW FedEx Standard Overnight®"    = "MW FedEx Standard Overnight"
    "MW FedEx 2Day® A.M."             = "MW FedEx 2Day A.M."
    "MW FedEx 2Day®"                  = "MW FedEx 2Day"
    "MW FedEx Express Saver®"         = "MW FedEx Express Saver"
    "FedEx First Overnight® Freight"  = "FedEx First Overnight Freight"
    "FedEx 1Day® Freight"             = "FedEx 1Day Freight"
    "FedEx 2Day® Freight"             = "FedEx 2Day Freight"
    "FedEx 3Day® Freight"             = "FedEx 3Day Freight"
}

foreach ($ws in $wb.Worksheets) {
    $oldName = $ws.Name
    if ($renames.ContainsKey($oldName)) {
        $ws.Name = $renames[$oldName]
    }
}

# Resize / reposition the workbook window to match the author's last
# on-screen layout.
$win = $excel.Windows.Item(1)
$win.Left = 28680
$win.Top = -1560
$win.Width = 29040
$win.Height = 15840
$win.TabRatio = 747

# Move the active tab off "FHD Ground" and onto "Intra-Hawaii Standard
# List Rate" - this flips tabSelected between the two sheetViews and
# updates workbookView/@activeTab.
$target = $wb.Worksheets.Item("Intra-Hawaii Standard List Rate")
$target.Activate()

# Scroll the tab strip so the newly active sheet's group is in view.
$win.ScrollWorkbookTabs(1, 5)

$wb.Save()
